$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I45").Value = 0.535526637811788
$ws.Range("H46").Value = 0.5766911554241068
$ws.Range("G47").Value = 0.6272238950261231
$ws.Range("F48").Value = 0.6666911554241067
$ws.Range("E49").Value = 0.6966911554241066
$ws.Range("D50").Value = 0.4271648845785767
$ws.Range("C51").Value = 0.4775315349050862
$ws.Range("B52").Value = 0.32386998960715
